$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# F5: "Endpointovi za pohranu podataka historiju podataka ..." ->
#     "Endpointovi za pohranu historije podataka ..."
$ws.Range("F5").Value = "Endpointovi za pohranu historije podataka zabilježenih sa senzora i deploy web api-a i baze podataka na fitov server."

# F9: insert ", korekcija u SQL bazi, dodavanje endpointova na WEB API" before
#     ", dodavanje Light komponente ..."
$ws.Range("F9").Value = "Dodavanje koda u Arduino IDE za upravljanje light senzorom, izmjene u Realtime DB na Firebase, korekcija u SQL bazi, dodavanje endpointova na WEB API, dodavanje Light komponente u Angularu i testiranje komunikacije na frontendu."

# Scroll the view so column D becomes the left-most visible column, and move
# the selection to F12 (matches the saved sheetView/selection state).
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F12").Select()
